# Weekly update: insert a fresh "latest" Choclo price record at the top of the
# data block (row 68) and push the existing history down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 68 - everything currently at row 68
# (and below, down to 121) shifts down to row 69 (through 122).
$ws.Rows("68:68").Insert()

# Populate the newly inserted row 68 with the new weekly record.
$ws.Range("A68").Value = 2
$ws.Range("B68").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C68").Value = "Coquimbo"
$ws.Range("D68").Value = 44658
$ws.Range("E68").Value = 4
$ws.Range("F68").Value = 100112024
$ws.Range("G68").Value = "Choclo"
$ws.Range("H68").Value = "Choclero"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 20000
$ws.Range("K68").Value = 200
$ws.Range("L68").Value = 230
$ws.Range("M68").Value = 215
$ws.Range("N68").Value = "$/unidad"
$ws.Range("O68").Value = "Provincia de Limarí"
$ws.Range("P68").Value = 215
$ws.Range("Q68").Value = 1
$ws.Range("R68").Value = "Hortaliza"
